# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks numeric must be forced to Text format first,
# otherwise Excel silently reinterprets them as floating-point numbers
# (losing the exact decimal text and changing the cell type).
$numericLookingCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D17",
    "D19",
    "D20",
    "D22",
    "D23",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D44",
    "D45",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Write the numeric-looking values as text.
$ws.Range("D5").Value = "302.79"
$ws.Range("D6").Value = "94.42"
$ws.Range("D7").Value = "0.501"
$ws.Range("D9").Value = "0.493"
$ws.Range("D10").Value = "34.09"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D12").Value = "18.73"
$ws.Range("D14").Value = "6.70"
$ws.Range("D17").Value = "0.791"
$ws.Range("D19").Value = "12.06"
$ws.Range("D20").Value = "6.19"
$ws.Range("D22").Value = "67.89"
$ws.Range("D23").Value = "235.53"
$ws.Range("D26").Value = "2.42"
$ws.Range("D27").Value = "24.46"
$ws.Range("D28").Value = "2.21"
$ws.Range("D29").Value = "9.12"
$ws.Range("D30").Value = "31.29"
$ws.Range("D33").Value = "4.98"
$ws.Range("D34").Value = "17.23"
$ws.Range("D35").Value = "4.38"
$ws.Range("D36").Value = "2.32"
$ws.Range("D37").Value = "1.81"
$ws.Range("D38").Value = "124.04"
$ws.Range("D40").Value = "2.73"
$ws.Range("D41").Value = "22.03"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D45").Value = "10.06"
$ws.Range("D47").Value = "2.72"
$ws.Range("D48").Value = "2.88"
$ws.Range("D50").Value = "52.54"
$ws.Range("D51").Value = "71.78"

# Restore default (General) styling now that the literal text is locked in,
# so the cell style matches the original unstyled cells.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining text/label/link/percentage cells - plain string values are fine.
$ws.Range("D2").Value = "42.940.98"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.330.23"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -3.46%  "
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "2.686.82"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.306.55"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "42.869.51"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("E19").Value = "  -4.79%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "0.0₃0886"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  +7.00%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  -6.43%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +7.99%  "
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E38").Value = "  -24.81%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E41").Value = "  +21.03%  "
$ws.Range("E42").Value = "  -1.43%  "
$ws.Range("D43").Value = "1.928.60"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  -6.18%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "2.556.56"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("E51").Value = "  -0.70%  "

